$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 600 (shifts old rows 600:641 down to 601:642)
$ws.Rows("600:600").Insert()

# Populate the newly inserted row with the new daily data point.
# Force the date column to stay as literal text (matching the other date
# cells in the column) rather than letting Excel auto-convert it to a
# date serial number.
$ws.Range("A600").NumberFormat = "@"
$ws.Range("A600").Value = "2026/01/11"
$ws.Range("A600").ClearFormats()

$ws.Range("B600").Value = "日"
$ws.Range("C600").Value = 8
$ws.Range("D600").Value = 20
